$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count and Wrong marking changes
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right total, Wrong total, and summary text
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "66 / 112"
